# Regenerate save_data to use K (Strike count) instead of Strike#.
# This updates column G ("K") values for each row of game data on the
# active worksheet to the freshly recomputed strike counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value, taken from the recomputed s_vals.
$kValues = [ordered]@{
    2  = 3
    3  = 0
    4  = 1
    6  = 1
    7  = 0
    8  = 2
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 3
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    29 = 3
    31 = 2
    32 = 1
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
